$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "开始日期" (start date) column (G). This shifts the old
# "结束日期" (end date) column H left into G, so the sheet ends up with a
# single date column (which we relabel below).
$ws.Columns("G").Delete()

# Relabel the remaining date column header as "到期日期" (due date).
$ws.Range("G1").Value = "到期日期"

# Update/append data rows 2-22. Columns: A=No, B=业务, C=QQ, D=VX, E=电话,
# F=网址, G=到期日期 (date).
$rows = @(
    @{ Row = 2;  B = "物流";  C = 418324114; D = "qqwee"; E = 18329953644; F = "https://raw.githubusercontent.com/Aishee001/ADUserdata/main/ADUserdata.xlsx"; G = "1/25/2026" }
    @{ Row = 3;  B = "账号";  C = 318324114;                E = 28329953645;                                                                                G = "1/25/2026" }
    @{ Row = 4;  B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 5;  B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 6;  B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "10/25/2025" }
    @{ Row = 7;  B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 8;  B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 9;  B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 10; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "12/20/2025" }
    @{ Row = 11; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 12; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 13; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 14; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/20/2025" }
    @{ Row = 15; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 16; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "11/25/2025" }
    @{ Row = 17; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 18; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 19; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "12/12/2025" }
    @{ Row = 20; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 21; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
    @{ Row = 22; B = "保险";  C = 218324114;                E = 38329953646;                                                                                G = "1/25/2026" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    if ($r.ContainsKey("D")) {
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    if ($r.ContainsKey("F")) {
        $ws.Cells.Item($r.Row, 6).Value = $r.F
    }
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}

# Match the narrower selection recorded in the saved workbook (single cell
# instead of a two-cell range).
[void]$ws.Range("J9").Select()
